$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update country name ordering (rank swaps due to refreshed totals) ---
$ws.Range("A47").Value = "Polonia"
$ws.Range("A48").Value = "Singapur"
$ws.Range("A71").Value = "Austria"
$ws.Range("A72").Value = "Australia"
$ws.Range("A199").Value = "Curazao"
$ws.Range("A200").Value = "Guam"
$ws.Range("A213").Value = "Montserrat"
$ws.Range("A214").Value = "Islas Malvinas"

# --- Update statistic columns (B..H) with refreshed figures ---
# Row 6
$ws.Range("B6").Value = 2530490
$ws.Range("C6").Value = 5268
$ws.Range("D6").Value = 1809702
$ws.Range("E6").Value = 671618
$ws.Range("G6").Value = 36
$ws.Range("H6").Value = 49170
# Row 7
$ws.Range("B7").Value = 917884
$ws.Range("C7").Value = 5061
$ws.Range("D7").Value = 729411
$ws.Range("E7").Value = 172856
$ws.Range("G7").Value = 119
$ws.Range("H7").Value = 15617
# Row 25
$ws.Range("B25").Value = 157918
$ws.Range("C25").Value = 4351
$ws.Range("D25").Value = 72209
$ws.Range("E25").Value = 83109
$ws.Range("G25").Value = 159
$ws.Range("H25").Value = 2600
# Row 34
$ws.Range("B34").Value = 89719
$ws.Range("C34").Value = 1847
$ws.Range("D34").Value = 47430
$ws.Range("E34").Value = 40245
$ws.Range("G34").Value = 33
$ws.Range("H34").Value = 2044
# Row 47
$ws.Range("B47").Value = 56090
$ws.Range("C47").Value = 771
$ws.Range("D47").Value = 38853
$ws.Range("E47").Value = 15368
$ws.Range("G47").Value = 11
$ws.Range("H47").Value = 1869
# Row 48
$ws.Range("B48").Value = 55580
$ws.Range("D48").Value = 51049
$ws.Range("E48").Value = 4504
$ws.Range("H48").Value = 27
# Row 71
$ws.Range("B71").Value = 23179
$ws.Range("C71").Value = 303
$ws.Range("D71").Value = 20627
$ws.Range("E71").Value = 1824
$ws.Range("G71").Value = 3
$ws.Range("H71").Value = 728
# Row 72
$ws.Range("B72").Value = 23035
$ws.Range("C72").Value = 292
$ws.Range("D72").Value = 13355
$ws.Range("E72").Value = 9301
$ws.Range("G72").Value = 4
$ws.Range("H72").Value = 379
# Row 73
$ws.Range("B73").Value = 22619
$ws.Range("C73").Value = 305
$ws.Range("D73").Value = 10618
$ws.Range("E73").Value = 11398
$ws.Range("G73").Value = 8
$ws.Range("H73").Value = 603
# Row 131
$ws.Range("B131").Value = 2184
$ws.Range("C131").Value = 7
$ws.Range("E131").Value = 145
# Row 147
$ws.Range("B147").Value = 1315
$ws.Range("C147").Value = 7
$ws.Range("E147").Value = 205
# Row 199
$ws.Range("B199").Value = 33
$ws.Range("C199").Value = 1
$ws.Range("D199").Value = 31
$ws.Range("E199").Value = 1
# Row 200
$ws.Range("D200").Value = 0
$ws.Range("E200").Value = 31
# Row 213
$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1
# Row 214
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

# --- Update "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 15 de Agosto de 2020 a las 10:38"
